# Merge the three runs of the "This program uses an API..." paragraph
# (on the "TextBox 9" shape of slide 1) into a single run, matching the
# author's edit that removed the mid-sentence run breaks around
# "associated endpoint".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by name so the script is resilient to shape ordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 9") {
        $shp = $s.Shapes.Item($i)
        break
    }
}

$tr = $shp.TextFrame.TextRange
$full = $tr.Text

$apos = [char]0x2019
$marker = "player$($apos)s name"

$startMarkerText = "This program uses an API"
$endMarkerText = "has played. "

$idxStart = $full.IndexOf($startMarkerText)
$idxEnd = $full.IndexOf($endMarkerText) + $endMarkerText.Length
$len = $idxEnd - $idxStart

$target = $shp.TextFrame.TextRange.Characters($idxStart + 1, $len)

$newText = "This program uses an API to display the statistics of a player based on user input. When the user enters a $marker, the program locates that player within the associated endpoint and returns the statistics for each season the player has played. "

$target.Text = $newText
